# CERS_QTR_FIN.xlsx quarterly update
# Inserts two new leftmost quarter columns (D, E) into the CERS sheet,
# shifting the existing quarterly data (old D:K) right to F:M, and
# populates the two new columns with the newest two quarters of data
# (period ending 2018-12-31 and 2018-09-30) for the Income Statement,
# Balance Sheet and Cash Flow Statement blocks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert two new columns at D (pushes existing D:K -> F:M)
# ---------------------------------------------------------------------
$ws.Range("D1:E1").EntireColumn.Insert()

# ---------------------------------------------------------------------
# 2. New data for the two inserted columns (D = newest quarter,
#    E = second-newest quarter), keyed by row number.
#    Rows holding a period-ending date use the date number format;
#    all other populated rows use the "#,##0" number format; rows
#    that are blank separators in the original data still receive the
#    "#,##0" number format on the new cells (no value).
# ---------------------------------------------------------------------
$dateRows = @(7, 38, 80)
$blankRows = @(11, 16, 19, 39, 40, 55, 56, 67, 82, 90, 95)

$rowData = @{
    7 = @(43465, 43373)
    8 = @(20200, 19300)
    9 = @(8400, 8100)
    10 = @(11800, 11200)
    12 = @(12400, 10800)
    13 = @(0, 0)
    14 = @(0, "NA")
    15 = @("NA", "NA")
    17 = @(35700, 32900)
    18 = @(-15500, -13600)
    20 = @(400, 500)
    21 = @(-14700, -12800)
    22 = @(1100, 1100)
    23 = @(-16100, -14100)
    24 = @(100, 100)
    25 = @(0, 0)
    26 = @(-16200, -14200)
    27 = @(-16200, -14200)
    28 = @(0, 0)
    29 = @(0, 0)
    30 = @(0, 0)
    31 = @(0, 0)
    32 = @(-400, -500)
    33 = @(-16200, -14200)
    34 = @(0, 0)
    35 = @(-16200, -14200)
    38 = @(43465, 43373)
    41 = @(28900, 22300)
    42 = @(88700, 96700)
    43 = @(8800, 10500)
    44 = @(13500, 13300)
    45 = @(7000, 7400)
    46 = @(146900, 150200)
    47 = @(0, 0)
    48 = @(8100, 3100)
    49 = @(1700, 1700)
    50 = @(0, 0)
    51 = @(0, 0)
    52 = @(6800, 6700)
    53 = @(0, 0)
    54 = @(163500, 161700)
    57 = @(18600, 13000)
    58 = @(7900, 5700)
    59 = @(26200, 24400)
    60 = @(52700, 43200)
    61 = @(22000, 24100)
    62 = @(4300, 2500)
    63 = @(0, 0)
    64 = @(0, 0)
    65 = @(0, 0)
    66 = @(78900, 69900)
    68 = @(0, 0)
    69 = @(0, 0)
    70 = @(0, 0)
    71 = @(0, 0)
    72 = @(-778900, -762700)
    73 = @(0, 0)
    74 = @(0, 0)
    75 = @(0, 0)
    76 = @(84500, 91800)
    77 = @(0, 0)
    80 = @(43465, 43373)
    81 = @(-16200, -14200)
    83 = @(300, 300)
    84 = @(0, 0)
    85 = @(0, 0)
    86 = @(0, 0)
    87 = @(0, 0)
    88 = @(0, 0)
    89 = @(-7100, -4900)
    91 = @(-500, -500)
    92 = @(0, 0)
    93 = @(0, 0)
    94 = @(7200, 500)
    96 = @(0, 0)
    97 = @(0, 0)
    98 = @(0, 0)
    99 = @(0, 0)
    100 = @(6400, 12400)
    101 = @(0, 0)
    102 = @(6500, 8000)
}

foreach ($row in $rowData.Keys) {
    $vals = $rowData[$row]
    $dCell = $ws.Cells.Item($row, 4)
    $eCell = $ws.Cells.Item($row, 5)

    if ($dateRows -contains $row) {
        $fmt = "[$-409]d\-mmm\-yy;@"
    } else {
        $fmt = "#,##0"
    }

    $dCell.NumberFormat = $fmt
    $eCell.NumberFormat = $fmt
    $dCell.Value2 = $vals[0]
    $eCell.Value2 = $vals[1]
}

foreach ($row in $blankRows) {
    $dCell = $ws.Cells.Item($row, 4)
    $eCell = $ws.Cells.Item($row, 5)
    $dCell.NumberFormat = "#,##0"
    $eCell.NumberFormat = "#,##0"
}

# ---------------------------------------------------------------------
# 3. Re-fit column widths to match the refreshed, wider data set.
# ---------------------------------------------------------------------
$ws.Columns.AutoFit()
